$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.164.46"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "1.562.59"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.14"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3799"
$ws.Range("E7").Value = "  +3.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3289"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.69"
$ws.Range("E9").Value = "  -9.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.141"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07372"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.01"
$ws.Range("E13").Value = "  -3.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.841"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.884"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "1.565.80"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001095"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06634"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.63"
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.471"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.14"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.76"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "22.176.25"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.265"
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.544"
$ws.Range("E26").Value = "  -3.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.44"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.11"
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.870"
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").Value = "1.733.75"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.57"
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.128"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.038"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.876"
$ws.Range("E34").Value = "  -5.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.373"
$ws.Range("E35").Value = "  -4.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08229"
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.303"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02309"
$ws.Range("E38").Value = "  -6.70%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2141"
$ws.Range("E40").Value = "  -4.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.233"
$ws.Range("E41").Value = "  -4.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.08"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5991"
$ws.Range("E44").Value = "  -4.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.74"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5800"
$ws.Range("E47").Value = "  -5.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.991"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.20"
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.173"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06993"
$ws.Range("E51").Value = "  -3.30%  "
